# Hindi-Imposition-Analysis.xlsx -- "Add files via upload" re-upload edit
#
# The underlying sentiment values for Neutral / Positive / Negative got
# shuffled around in the source table (Sheet1!B3:B5). Re-enter them in
# their new order and move the active selection to B7 (matching the
# freshly-saved workbook's cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B3 (Neutral sentiment):  0.0942  -> 0.8692
$ws.Range("B3").Value = 0.86916000000000004

# B4 (Positive sentiment): 0.8692  -> 0.0366
$ws.Range("B4").Value = 0.036639999999999999

# B5 (Negative sentiment): 0.0366  -> 0.0942
$ws.Range("B5").Value = 0.094180000000000097

# Move the selection/active cell to B7, as in the saved workbook.
$ws.Range("B7").Select()
